# Generate Report for Handback
# Update the handoff/handback timestamps for the 1a11b725... (zh-cn / de-de) rows
# on the "zh-cn" and "de-de" worksheets (row 3, columns E and H).

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-21 02:36:05"
$wsZhCn.Range("H3").Value = "2016-03-21 02:36:27"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-21 02:36:08"
$wsDeDe.Range("H3").Value = "2016-03-21 02:36:33"
